# #435 removed slug column
#
# The "Products" sheet has a "Slug" header column (with values like
# "Cherry Mobile" underneath) that is no longer needed. Remove that whole
# column so every column to its right shifts left by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")
$ws.Activate()

# Locate the "Slug" header in row 1 and delete its entire column.
$headerRow = $ws.Range("A1:Z1")
$slugHeader = $headerRow.Find("Slug")
if ($slugHeader -ne $null) {
    $slugHeader.EntireColumn.Delete()
}

# Mirror the author's resulting selection on the Products sheet.
$ws.Range("L4").Select()
